$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the data ranges so numeric-looking strings
# (prices, percentages, hour values) are preserved as text, matching the
# workbook author's original inline-string cells instead of being coerced
# into Excel numbers/percentages.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Updated symbol list values (price, volume(1h), hour) scraped on
# Wed Jan  4 21:08:39 UTC 2023
$ws.Range("D2").Value = "256.43"
$ws.Range("E2").Value = "4.50%"
$ws.Range("G2").Value = "21"
$ws.Range("D3").Value = "27.54"
$ws.Range("E3").Value = "-2.42%"
$ws.Range("G3").Value = "21"
$ws.Range("D4").Value = "5.213"
$ws.Range("E4").Value = "-0.83%"
$ws.Range("G4").Value = "21"
$ws.Range("D5").Value = "0.05922"
$ws.Range("E5").Value = "3.86%"
$ws.Range("G5").Value = "21"
$ws.Range("D6").Value = "6.668"
$ws.Range("E6").Value = "0.42%"
$ws.Range("G6").Value = "21"
$ws.Range("D7").Value = "0.8666"
$ws.Range("E7").Value = "1.83%"
$ws.Range("G7").Value = "21"
$ws.Range("E8").Value = "14.69%"
$ws.Range("G8").Value = "21"
$ws.Range("D9").Value = "0.1419"
$ws.Range("E9").Value = "1.26%"
$ws.Range("G9").Value = "21"
$ws.Range("D10").Value = "0.07185"
$ws.Range("E10").Value = "1.36%"
$ws.Range("G10").Value = "21"
$ws.Range("E11").Value = "3.05%"
$ws.Range("G11").Value = "21"
$ws.Range("D12").Value = "0.09218"
$ws.Range("E12").Value = "-0.12%"
$ws.Range("G12").Value = "21"
$ws.Range("D13").Value = "0.001543"
$ws.Range("E13").Value = "1.16%"
$ws.Range("G13").Value = "21"
$ws.Range("D14").Value = "0.0006062"
$ws.Range("E14").Value = "-93.96%"
$ws.Range("G14").Value = "21"
$ws.Range("D15").Value = "0.005722"
$ws.Range("E15").Value = "-4.44%"
$ws.Range("G15").Value = "21"
$ws.Range("D16").Value = "3.479"
$ws.Range("E16").Value = "-0.31%"
$ws.Range("G16").Value = "21"
$ws.Range("D17").Value = "3.266"
$ws.Range("E17").Value = "1.87%"
$ws.Range("G17").Value = "21"
$ws.Range("D18").Value = "2.224"
$ws.Range("E18").Value = "1.62%"
$ws.Range("G18").Value = "21"
$ws.Range("D19").Value = "0.3150"
$ws.Range("E19").Value = "-0.54%"
$ws.Range("G19").Value = "21"
$ws.Range("D20").Value = "0.03610"
$ws.Range("E20").Value = "9.11%"
$ws.Range("G20").Value = "21"
$ws.Range("D21").Value = "0.1307"
$ws.Range("E21").Value = "2.77%"
$ws.Range("G21").Value = "21"
$ws.Range("D22").Value = "3.524"
$ws.Range("E22").Value = "-0.23%"
$ws.Range("G22").Value = "21"
$ws.Range("D23").Value = "0.04180"
$ws.Range("E23").Value = "2.88%"
$ws.Range("G23").Value = "21"
$ws.Range("D24").Value = "0.1399"
$ws.Range("E24").Value = "1.49%"
$ws.Range("G24").Value = "21"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").Value = "-0.33%"
$ws.Range("G25").Value = "21"
$ws.Range("D26").Value = "0.004517"
$ws.Range("E26").Value = "8.84%"
$ws.Range("G26").Value = "21"
$ws.Range("E27").Value = "0.09%"
$ws.Range("G27").Value = "21"
$ws.Range("E28").Value = "33.80%"
$ws.Range("G28").Value = "21"
$ws.Range("G29").Value = "21"
$ws.Range("G30").Value = "21"
$ws.Range("G31").Value = "21"
$ws.Range("G32").Value = "21"
$ws.Range("G33").Value = "21"
$ws.Range("G34").Value = "21"
$ws.Range("G35").Value = "21"
$ws.Range("G36").Value = "21"
$ws.Range("G37").Value = "21"
$ws.Range("G38").Value = "21"
$ws.Range("G39").Value = "21"
$ws.Range("D40").Value = "0.03821"
$ws.Range("E40").Value = "0.64%"
$ws.Range("G40").Value = "21"
$ws.Range("D41").Value = "0.005504"
$ws.Range("E41").Value = "6.70%"
$ws.Range("G41").Value = "21"
$ws.Range("D42").Value = "0.1102"
$ws.Range("E42").Value = "3.49%"
$ws.Range("G42").Value = "21"
$ws.Range("E43").Value = "-13.56%"
$ws.Range("G43").Value = "21"
$ws.Range("D44").Value = "0.01069"
$ws.Range("E44").Value = "3.63%"
$ws.Range("G44").Value = "21"
$ws.Range("D45").Value = "0.00005432"
$ws.Range("E45").Value = "3.13%"
$ws.Range("G45").Value = "21"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.09%"
$ws.Range("G46").Value = "21"
$ws.Range("D47").Value = "0.1091"
$ws.Range("E47").Value = "3.94%"
$ws.Range("G47").Value = "21"
$ws.Range("D48").Value = "0.002172"
$ws.Range("E48").Value = "-4.32%"
$ws.Range("G48").Value = "21"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "0.09%"
$ws.Range("G49").Value = "21"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "21"
$ws.Range("G51").Value = "21"
